# --- Spice.xlsx: add "Work Order Item" / "Quantity" / "Lot" columns to the
#     consumption template, and introduce a second "lookup" sheet. ---

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "main"

# New "lookup" sheet, placed right after "main".
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "lookup"

# Shift the existing headers one column to the right to make room for the
# new "Work Order Item" column, then fill in the new headers.
$ws.Columns.Item(1).Insert()
$ws.Range("A2").Value = "Work Order Item"
$ws.Range("F2").Value = "Quantity"
$ws.Range("E2").Value = "Lot"

# Match the yellow header highlight used by the rest of the row.
$ws.Range("E2:F2").Interior.Color = 65535

# Column widths (best-fit sizing for the new header text).
$ws.Columns.Item(1).ColumnWidth = 14.333333333333334
$ws.Columns.Item(3).ColumnWidth = 9.5
$ws.Columns.Item(5).ColumnWidth = 24.666666666666668
$ws.Columns.Item(6).ColumnWidth = 18

$ws2.Columns.Item(1).ColumnWidth = 15.0
$ws2.Columns.Item(2).ColumnWidth = 25.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 8.833333333333334
$ws2.Columns.Item(4).ColumnWidth = 15.5
$ws2.Columns.Item(5).ColumnWidth = 3.5
$ws2.Columns.Item(6).ColumnWidth = 3.5

# Selections left behind in the source workbook.
$ws2.Rows("1:19").Select() | Out-Null
$ws.Activate()
$ws.Range("F6").Select() | Out-Null
